$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=258; B="Ilyes Boughanmi"; D=7; E=6; F=0; G=$null; H=0 },
    @{ Row=259; B="Omar Benyounes"; D=6; E=6; F=0; G=$null; H=3 },
    @{ Row=260; B="Amir Etien"; D=8; E=10; F=0; G=$null; H=4 },
    @{ Row=261; B="Romain Thunet"; D=8; E=8; F=7; G="Béquille "; H=1 },
    @{ Row=262; B="Rayane Chayebi"; D=9; E=10; F=8; G="Adducteur "; H=5 },
    @{ Row=263; B="Naim Ighbane"; D=7; E=0; F=0; G=$null; H=4 },
    @{ Row=264; B="Jeremie Laurent"; D=6; E=6; F=0; G=$null; H=7 },
    @{ Row=265; B="Ilan Ihaddadene"; D=7; E=8; F=0; G=$null; H=8 },
    @{ Row=266; B="Naim Dhib"; D=6; E=7; F=3; G="Quadri courbature"; H=5 },
    @{ Row=267; B="Sofiane Belle"; D=6; E=4; F=0; G=$null; H=5 },
    @{ Row=268; B="Amir Kherrab"; D=4; E=6; F=6; G="Semelle "; H=5 },
    @{ Row=269; B="Wael Fareh"; D=6; E=3; F=1; G="Genou "; H=6 },
    @{ Row=270; B="Hedi Nasri"; D=7; E=7; F=3; G="Adducteur "; H=7 },
    @{ Row=271; B="Yoan Zouma"; D=4; E=5; F=3; G="Cheville"; H=5 },
    @{ Row=272; B="Yanis Berrached"; D=5; E=10; F=0; G=$null; H=0 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $fullRange = "A" + $rowNum + ":I" + $rowNum
    if ($r.G -eq $null) {
        # template row with empty "Localisation douleur" (style s="2")
        $ws.Range("A250:I250").Copy($ws.Range($fullRange))
    } else {
        # template row with text in "Localisation douleur" (style s="1")
        $ws.Range("A257:I257").Copy($ws.Range($fullRange))
    }
    $excel.CutCopyMode = 0

    $ws.Range("A" + $rowNum).Value = 45903
    $ws.Range("B" + $rowNum).Value = $r.B
    $ws.Range("C" + $rowNum).Value = 70
    $ws.Range("D" + $rowNum).Value = $r.D
    $ws.Range("E" + $rowNum).Value = $r.E
    $ws.Range("F" + $rowNum).Value = $r.F
    if ($r.G -ne $null) {
        $ws.Range("G" + $rowNum).Value = $r.G
    }
    $ws.Range("H" + $rowNum).Value = $r.H
    $ws.Range("I" + $rowNum).Formula = "=C" + $rowNum + "*D" + $rowNum
}

$ws.Range("K266").Select() | Out-Null

Write-Host "done"